$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 30499.666
$ws.Range("J3").Value = 30499.666
$ws.Range("L3").Value = 30499.666
$ws.Range("N3").Value = -30727.666

$ws.Range("H53").Value = 633.0909
$ws.Range("I53").Value = 696.125
$ws.Range("K53").Value = 696.125
$ws.Range("M53").Value = -59.125

$ws.Range("H82").Value = 647.8
$ws.Range("I82").Value = 647.8
$ws.Range("K82").Value = 1943.4
$ws.Range("M82").Value = -1537.4

$ws.Range("H85").Value = 647.8
$ws.Range("I85").Value = 647.8
$ws.Range("K85").Value = 1943.4
$ws.Range("M85").Value = -539.3999999999999

$ws.Range("H102").Value = 30499.666
$ws.Range("J102").Value = 30499.666
$ws.Range("L102").Value = 30499.666
$ws.Range("N102").Value = -36989.666

$ws.Range("H104").Value = 1150
$ws.Range("I104").Value = 1150
$ws.Range("K104").Value = 3450
$ws.Range("M104").Value = -1703

$ws.Range("H113").Value = 1500
$ws.Range("I113").Value = 1500
$ws.Range("K113").Value = 1500
$ws.Range("M113").Value = 1754

$ws.Range("H125").Value = 3954.8
$ws.Range("I125").Value = 3693.5
$ws.Range("K125").Value = 33241.5
$ws.Range("M125").Value = -30781.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2420.5625
$ws.Range("I2").Value = 2364.5386
$ws.Range("J2").Value = 2663.3333
$ws.Range("K2").Value = 2364.5386
$ws.Range("L2").Value = 2663.3333
$ws.Range("M2").Value = -2251.5386
$ws.Range("N2").Value = -2889.3333

$ws.Range("H45").Value = 2996.4285
$ws.Range("I45").Value = 2795
$ws.Range("J45").Value = 3500
$ws.Range("K45").Value = 2795
$ws.Range("L45").Value = 3500
$ws.Range("M45").Value = -2418
$ws.Range("N45").Value = -4254

$ws.Range("H74").Value = 6182.909
$ws.Range("I74").Value = 5099.6
$ws.Range("J74").Value = 7085.6665
$ws.Range("K74").Value = 5099.6
$ws.Range("L74").Value = 7085.6665
$ws.Range("M74").Value = -4225.6
$ws.Range("N74").Value = -8833.666499999999

$ws.Range("H77").Value = 6182.909
$ws.Range("I77").Value = 5099.6
$ws.Range("J77").Value = 7085.6665
$ws.Range("K77").Value = 25498
$ws.Range("L77").Value = 35428.3325
$ws.Range("M77").Value = -21130
$ws.Range("N77").Value = -44164.3325

$ws.Range("H102").Value = 6134.25
$ws.Range("I102").Value = 4613.8
$ws.Range("K102").Value = 4613.8
$ws.Range("M102").Value = -2991.8

$ws.Range("H116").Value = 2420.5625
$ws.Range("I116").Value = 2364.5386
$ws.Range("J116").Value = 2663.3333
$ws.Range("K116").Value = 2364.5386
$ws.Range("L116").Value = 2663.3333
$ws.Range("M116").Value = -70.53859999999986
$ws.Range("N116").Value = -7251.3333

$ws.Range("H122").Value = 1766.6666
$ws.Range("I122").Value = 1766.6666
$ws.Range("K122").Value = 5299.9998
$ws.Range("M122").Value = -2849.9998

$ws.Range("H124").Value = 75000
$ws.Range("J124").Value = 75000
$ws.Range("L124").Value = 75000
$ws.Range("N124").Value = -84820

$ws.Range("H125").Value = 70333.336
$ws.Range("J125").Value = 70333.336
$ws.Range("L125").Value = 70333.336
$ws.Range("N125").Value = -80173.336

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2420.5625
$ws.Range("I3").Value = 2364.5386
$ws.Range("J3").Value = 2663.3333
$ws.Range("K3").Value = 2364.5386
$ws.Range("L3").Value = 2663.3333
$ws.Range("M3").Value = -2250.5386
$ws.Range("N3").Value = -2891.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2497.5
$ws.Range("J2").Value = 2497.5
$ws.Range("L2").Value = 2497.5
$ws.Range("N2").Value = -2723.5

$ws.Range("H99").Value = 3890.7334
$ws.Range("J99").Value = 4604.3335
$ws.Range("L99").Value = 4604.3335
$ws.Range("N99").Value = -7600.3335

$ws.Range("H124").Value = 30000
$ws.Range("I124").Value = 30000
$ws.Range("K124").Value = 30000
$ws.Range("M124").Value = -27545

$ws.Range("H126").Value = 3890.7334
$ws.Range("J126").Value = 4604.3335
$ws.Range("L126").Value = 13813.0005
$ws.Range("N126").Value = -18753.0005

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()

$ws.Range("H63").Value = 1949.5
$ws.Range("I63").Value = 1949.5
$ws.Range("K63").Value = 5848.5
$ws.Range("M63").Value = -5099.5

$ws.Range("H66").Value = 1949.5
$ws.Range("I66").Value = 1949.5
$ws.Range("K66").Value = 17545.5
$ws.Range("M66").Value = -13801.5

$ws.Range("H102").Value = 5000
$ws.Range("I102").Value = 5000
$ws.Range("K102").Value = 15000
$ws.Range("M102").Value = -12566

$ws.Range("H104").Value = 12500
$ws.Range("I104").Value = 6666.6665
$ws.Range("J104").Value = 16000
$ws.Range("K104").Value = 19999.9995
$ws.Range("L104").Value = 48000
$ws.Range("M104").Value = -17378.9995
$ws.Range("N104").Value = -53242

$ws.Range("H114").Value = 947
$ws.Range("I114").Value = 1019.3333
$ws.Range("J114").Value = 874.6667
$ws.Range("K114").Value = 3057.9999
$ws.Range("L114").Value = 2624.0001
$ws.Range("M114").Value = 196.0001000000002
$ws.Range("N114").Value = -9132.000100000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

$ws.Range("H80").Value = 2247.6
$ws.Range("J80").Value = 2247.6
$ws.Range("L80").Value = 2247.6
$ws.Range("N80").Value = -4243.6

$ws.Range("H83").Value = 2247.6
$ws.Range("J83").Value = 2247.6
$ws.Range("L83").Value = 11238
$ws.Range("N83").Value = -21222

$ws.Range("H95").Value = 22524.5
$ws.Range("J95").Value = 22524.5
$ws.Range("L95").Value = 22524.5
$ws.Range("N95").Value = -28016.5

$ws.Range("H99").Value = 6494.2
$ws.Range("I99").Value = 4117.75
$ws.Range("J99").Value = 16000
$ws.Range("K99").Value = 4117.75
$ws.Range("L99").Value = 16000
$ws.Range("M99").Value = -1871.75
$ws.Range("N99").Value = -20492

$ws.Range("H113").Value = 7396.5
$ws.Range("I113").Value = 4966.25
$ws.Range("J113").Value = 8611.625
$ws.Range("K113").Value = 4966.25
$ws.Range("L113").Value = 8611.625
$ws.Range("M113").Value = -2796.25
$ws.Range("N113").Value = -12951.625

$ws.Range("H122").Value = 4354.6665
$ws.Range("I122").Value = 3670.2
$ws.Range("K122").Value = 11010.6
$ws.Range("M122").Value = -8560.599999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5351.9165
$ws.Range("I7").Value = 5235.25
$ws.Range("J7").Value = 5585.25
$ws.Range("K7").Value = 5235.25
$ws.Range("L7").Value = 5585.25
$ws.Range("M7").Value = -5123.25
$ws.Range("N7").Value = -5809.25

$ws.Range("H40").Value = 9575
$ws.Range("J40").Value = 9433.333000000001
$ws.Range("L40").Value = 9433.333000000001
$ws.Range("N40").Value = -9705.333000000001

$ws.Range("H61").Value = 6697.4375
$ws.Range("I61").Value = 5715.375
$ws.Range("J61").Value = 7679.5
$ws.Range("K61").Value = 5715.375
$ws.Range("L61").Value = 7679.5
$ws.Range("M61").Value = -5513.375
$ws.Range("N61").Value = -8083.5

$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()

$ws.Range("H113").Value = 6697.4375
$ws.Range("I113").Value = 5715.375
$ws.Range("J113").Value = 7679.5
$ws.Range("K113").Value = 5715.375
$ws.Range("L113").Value = 7679.5
$ws.Range("M113").Value = -3545.375
$ws.Range("N113").Value = -12019.5

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H126").Value = 5351.9165
$ws.Range("I126").Value = 5235.25
$ws.Range("J126").Value = 5585.25
$ws.Range("K126").Value = 15705.75
$ws.Range("L126").Value = 16755.75
$ws.Range("M126").Value = -13235.75
$ws.Range("N126").Value = -21695.75

$ws.Range("H132").Value = 5962.2
$ws.Range("I132").Value = 4220.5
$ws.Range("J132").Value = 7123.3335
$ws.Range("K132").Value = 12661.5
$ws.Range("L132").Value = 21370.0005
$ws.Range("M132").Value = -10131.5
$ws.Range("N132").Value = -26430.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 185163.64
$ws.Range("I2").Value = 185163.64
$ws.Range("K2").Value = 185163.64
$ws.Range("M2").Value = -185051.64

$ws.Range("H4").Value = 163176.92
$ws.Range("I4").Value = 212050
$ws.Range("K4").Value = 212050
$ws.Range("M4").Value = -211937

$ws.Range("H100").Value = 683.75
$ws.Range("I100").Value = 683.75
$ws.Range("K100").Value = 1367.5
$ws.Range("M100").Value = -826.5

$ws.Range("H113").Value = 997.5714
$ws.Range("I113").Value = 1764
$ws.Range("K113").Value = 5292
$ws.Range("M113").Value = -3122

$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 2000
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 6000
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -3550
$ws.Range("N122").Value = -19900

$ws.Range("H126").Value = 3325.2
$ws.Range("I126").Value = 1357.0714
$ws.Range("K126").Value = 4071.2142
$ws.Range("M126").Value = -1601.2142

$ws.Range("H132").Value = 3097.4546
$ws.Range("I132").Value = 3097.4546
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 9292.363799999999
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").Value = -6762.363799999999
